$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 42605.886712962965
$ws.Cells.Item(4,2).Value = 46
$ws.Cells.Item(4,3).Value = 0
$ws.Cells.Item(4,4).Value = 0
$ws.Cells.Item(4,5).Value = 0
$ws.Cells.Item(4,6).Value = 0
$ws.Cells.Item(4,7).Value = 0
$ws.Cells.Item(4,8).Value = 0
$ws.Cells.Item(4,9).Value = 0
$ws.Cells.Item(4,10).Value = 0
$ws.Cells.Item(4,11).Value = 0
$ws.Cells.Item(4,12).Value = 0
$ws.Cells.Item(4,13).Value = 0
$ws.Cells.Item(4,14).Value = "Random"
